$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 345.45456
$ws.Range("I33").Value = 281.57144
$ws.Range("J33").Value = 457.25
$ws.Range("K33").Value = 281.57144
$ws.Range("L33").Value = 457.25
$ws.Range("M33").Value = -52.57144
$ws.Range("N33").Value = -915.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 350.05264
$ws.Range("I80").Value = 322
$ws.Range("J80").Value = 499.66666
$ws.Range("K80").Value = 966
$ws.Range("L80").Value = 1498.99998
$ws.Range("M80").Value = 32
$ws.Range("N80").Value = -3494.99998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 350.05264
$ws.Range("I83").Value = 322
$ws.Range("J83").Value = 499.66666
$ws.Range("K83").Value = 2898
$ws.Range("L83").Value = 4496.99994
$ws.Range("M83").Value = 2094
$ws.Range("N83").Value = -14480.99994

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 979.6
$ws.Range("I125").Value = 949
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 8541
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -6081
$ws.Range("N125").Value = -13920

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 46580
$ws.Range("J134").Value = 46580
$ws.Range("L134").Value = 46580
$ws.Range("N134").Value = -56720

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 794.3415
$ws.Range("I74").Value = 770.5143
$ws.Range("J74").Value = 933.3333
$ws.Range("K74").Value = 770.5143
$ws.Range("L74").Value = 933.3333
$ws.Range("M74").Value = 103.4857
$ws.Range("N74").Value = -2681.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 794.3415
$ws.Range("I77").Value = 770.5143
$ws.Range("J77").Value = 933.3333
$ws.Range("K77").Value = 3852.5715
$ws.Range("L77").Value = 4666.6665
$ws.Range("M77").Value = 515.4285
$ws.Range("N77").Value = -13402.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 13933175
$ws.Range("I132").Value = 16667762
$ws.Range("J132").Value = 3678476.5
$ws.Range("K132").Value = 50003286
$ws.Range("L132").Value = 11035429.5
$ws.Range("M132").Value = -50000756
$ws.Range("N132").Value = -11040489.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3471.6667
$ws.Range("I86").Value = 4388.2
$ws.Range("K86").Value = 4388.2
$ws.Range("M86").Value = -3265.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3471.6667
$ws.Range("I89").Value = 4388.2
$ws.Range("K89").Value = 21941
$ws.Range("M89").Value = -16325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2528522.8
$ws.Range("I134").Value = 853.5862
$ws.Range("K134").Value = 2560.7586
$ws.Range("M134").Value = -25.75860000000011

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 30306086
$ws.Range("I99").Value = 52634230
$ws.Range("J99").Value = 3600
$ws.Range("K99").Value = 52634230
$ws.Range("L99").Value = 3600
$ws.Range("M99").Value = -52632732
$ws.Range("N99").Value = -6596

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 30306086
$ws.Range("I126").Value = 52634230
$ws.Range("J126").Value = 3600
$ws.Range("K126").Value = 157902690
$ws.Range("L126").Value = 10800
$ws.Range("M126").Value = -157900220
$ws.Range("N126").Value = -15740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9260530
$ws.Range("I132").Value = 1194.0358
$ws.Range("J132").Value = 41668210
$ws.Range("K132").Value = 3582.1074
$ws.Range("L132").Value = 125004630
$ws.Range("M132").Value = -1052.1074
$ws.Range("N132").Value = -125009690

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 25588.535
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 34383.562
$ws.Range("K12").Value = 9
$ws.Range("L12").Value = 103150.686
$ws.Range("M12").Value = 164
$ws.Range("N12").Value = -103496.686

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 23843042
$ws.Range("I113").Value = 9259692
$ws.Range("J113").Value = 32593054
$ws.Range("K113").Value = 27779076
$ws.Range("L113").Value = 97779162
$ws.Range("M113").Value = -27776906
$ws.Range("N113").Value = -97783502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 10780300
$ws.Range("I122").Value = 69444700
$ws.Range("J122").Value = 5205.9185
$ws.Range("K122").Value = 625002300
$ws.Range("L122").Value = 46853.2665
$ws.Range("M122").Value = -624999850
$ws.Range("N122").Value = -51753.2665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 26322584
$ws.Range("I132").Value = 947
$ws.Range("J132").Value = 33341686
$ws.Range("K132").Value = 8523
$ws.Range("L132").Value = 300075174
$ws.Range("M132").Value = -5993
$ws.Range("N132").Value = -300080234

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 41669216
$ws.Range("I133").Value = 47621270
$ws.Range("J133").Value = 4833
$ws.Range("K133").Value = 142863810
$ws.Range("L133").Value = 14499
$ws.Range("M133").Value = -142858750
$ws.Range("N133").Value = -24619

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 17246084
$ws.Range("I122").Value = 21744580
$ws.Range("K122").Value = 65233740
$ws.Range("M122").Value = -65231290

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2264.5
$ws.Range("I126").Value = 1288.8889
$ws.Range("J126").Value = 3062.7273
$ws.Range("K126").Value = 3866.6667
$ws.Range("L126").Value = 9188.1819
$ws.Range("M126").Value = -1396.6667
$ws.Range("N126").Value = -14128.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7417.4
$ws.Range("I132").Value = 4784.88
$ws.Range("J132").Value = 20580
$ws.Range("K132").Value = 14354.64
$ws.Range("L132").Value = 61740
$ws.Range("M132").Value = -11824.64
$ws.Range("N132").Value = -66800

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1940.4286
$ws.Range("I7").Value = 1845.75
$ws.Range("J7").Value = 2066.6667
$ws.Range("K7").Value = 1845.75
$ws.Range("L7").Value = 2066.6667
$ws.Range("M7").Value = -1733.75
$ws.Range("N7").Value = -2290.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 57890
$ws.Range("J36").Value = 57890
$ws.Range("L36").Value = 57890
$ws.Range("N36").Value = -59014

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7943.1763
$ws.Range("I122").Value = 9542.615
$ws.Range("J122").Value = 2745
$ws.Range("K122").Value = 28627.845
$ws.Range("L122").Value = 8235
$ws.Range("M122").Value = -26177.845
$ws.Range("N122").Value = -13135

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1940.4286
$ws.Range("I126").Value = 1845.75
$ws.Range("J126").Value = 2066.6667
$ws.Range("K126").Value = 5537.25
$ws.Range("L126").Value = 6200.000100000001
$ws.Range("M126").Value = -3067.25
$ws.Range("N126").Value = -11140.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 23262718
$ws.Range("I132").Value = 32259964
$ws.Range("J132").Value = 19832.084
$ws.Range("K132").Value = 96779892
$ws.Range("L132").Value = 59496.25199999999
$ws.Range("M132").Value = -96777362
$ws.Range("N132").Value = -64556.25199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 17973.387
$ws.Range("I122").Value = 30823.47
$ws.Range("K122").Value = 92470.41
$ws.Range("M122").Value = -90020.41

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1275.2307
$ws.Range("I126").Value = 864.2222
$ws.Range("K126").Value = 2592.6666
$ws.Range("M126").Value = -122.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18876.188
$ws.Range("I132").Value = 20910.725
$ws.Range("J132").Value = 10894.538
$ws.Range("K132").Value = 62732.175
$ws.Range("L132").Value = 32683.614
$ws.Range("M132").Value = -60202.175
$ws.Range("N132").Value = -37743.614
